$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11
$t.Cell(1,1).Range.Text = '71 x 20' + $nl + '  2    0' + $nl + '  ----' + $nl + '7|    |' + $nl + '1|    |'
$t.Cell(1,2).Range.Text = '51 x 11' + $nl + '  1    1' + $nl + '  ----' + $nl + '5|    |' + $nl + '1|    |'
$t.Cell(1,3).Range.Text = '98 x 38' + $nl + '  3    8' + $nl + '  ----' + $nl + '9|    |' + $nl + '8|    |'
$t.Cell(2,1).Range.Text = '59 x 40' + $nl + '  4    0' + $nl + '  ----' + $nl + '5|    |' + $nl + '9|    |'
$t.Cell(2,2).Range.Text = '75 x 87' + $nl + '  8    7' + $nl + '  ----' + $nl + '7|    |' + $nl + '5|    |'
$t.Cell(2,3).Range.Text = '52 x 84' + $nl + '  8    4' + $nl + '  ----' + $nl + '5|    |' + $nl + '2|    |'
$t.Cell(3,1).Range.Text = '73 x 85' + $nl + '  8    5' + $nl + '  ----' + $nl + '7|    |' + $nl + '3|    |'
$t.Cell(3,2).Range.Text = '24 x 82' + $nl + '  8    2' + $nl + '  ----' + $nl + '2|    |' + $nl + '4|    |'
$t.Cell(3,3).Range.Text = '77 x 48' + $nl + '  4    8' + $nl + '  ----' + $nl + '7|    |' + $nl + '7|    |'
$t.Cell(4,1).Range.Text = '85 x 62' + $nl + '  6    2' + $nl + '  ----' + $nl + '8|    |' + $nl + '5|    |'
$t.Cell(4,2).Range.Text = '79 x 46' + $nl + '  4    6' + $nl + '  ----' + $nl + '7|    |' + $nl + '9|    |'
$t.Cell(4,3).Range.Text = '19 x 58' + $nl + '  5    8' + $nl + '  ----' + $nl + '1|    |' + $nl + '9|    |'
$t.Cell(5,1).Range.Text = '17 x 53' + $nl + '  5    3' + $nl + '  ----' + $nl + '1|    |' + $nl + '7|    |'
$t.Cell(5,2).Range.Text = '98 x 78' + $nl + '  7    8' + $nl + '  ----' + $nl + '9|    |' + $nl + '8|    |'
$t.Cell(5,3).Range.Text = '62 x 78' + $nl + '  7    8' + $nl + '  ----' + $nl + '6|    |' + $nl + '2|    |'